# Update the "Förändrad" (Changed) date column (C) for data rows 2-31
# from 45594 (2024-10-29) to 45595 (2024-10-30).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45594) {
        $cell.Value2 = 45595
    }
}
